$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New training-run results appended by the training pipeline (rows 10-13).
# Note: column F (Test R2) is intentionally left blank for these rows.
$rows = @(
    @{ Row = 10; A = "2024-11-19 18:25:43"; B = 0.9946636683083764918222869; C = 0.0104692156432963593554186; D = 0.0002421817431318169105709; E = 0.0155621895352748094293327; G = 0.0116096623221671607639749; H = 0.0001618751108549740084803; I = 0.0127230150064744501614555 },
    @{ Row = 11; A = "2024-11-19 19:02:54"; B = 0.9963495219251747858635326; C = 0.0073396607095128421627273; D = 0.0001656700731818863095608; E = 0.0128712887148834595474334; G = 0.0086048133651053897164029; H = 0.0000875497349478999813133; I = 0.0093568015340660063566824 },
    @{ Row = 12; A = "2024-11-19 19:44:13"; B = 0.9963513355189759224117552; C = 0.0073403674720669640635373; D = 0.0001656264300690109920675; E = 0.0128695932363463196035358; G = 0.0147264909429365691700964; H = 0.0002746907459648597216867; I = 0.0165737969688559798331351 },
    @{ Row = 13; A = "2024-11-19 19:48:06"; B = 0.9957472806456137970343434; C = 0.0088764321778146641067631; D = 0.0001930474429988412066617; E = 0.0138941513954196297492283; G = 0.0136749457998091208249525; H = 0.0002620707470811911164905; I = 0.0161885992933666414406169 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
}

# Reflect the selection state recorded after the edit.
$ws.Range("F15").Select()
